$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift data rows up by one: row 3 (old "Respondent ratio (%)") is removed,
#     rows 4-11 shift up to become rows 3-10; then remove the new row 10
#     (the second now-blank trailer row), leaving 9 rows total (A1:R9). ---
$ws.Rows("3:3").Delete()
$ws.Rows("10:10").Delete()

# --- Row 1: clear A1 text, and strip bold/border/center-align formatting from header row ---
$ws.Range("A1").Value = ""
$ws.Range("A1:R1").Style = "Normal"

# --- Row 3: corrected data values ---
$ws.Range("A3").Value = "Revisit count"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 13
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 37
$ws.Range("Q3").Value = 4

# --- Row 4: corrected data values ---
$ws.Range("A4").Value = "Fixation count"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 122
$ws.Range("E4").Value = 27
$ws.Range("F4").Value = 24
$ws.Range("G4").Value = 188
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 21
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 119
$ws.Range("Q4").Value = 6

# --- Row 5: corrected data values ---
$ws.Range("A5").Value = "Dwell time (ms)"
$ws.Range("B5").Value = 2385.63
$ws.Range("C5").Value = 700.59
$ws.Range("D5").Value = 35302.28
$ws.Range("E5").Value = 8708.2
$ws.Range("F5").Value = 6856.85
$ws.Range("G5").Value = 47517.33
$ws.Range("H5").Value = 1401.57
$ws.Range("I5").Value = 1751.72
$ws.Range("J5").Value = 350.36
$ws.Range("K5").Value = 1067.48
$ws.Range("L5").Value = 1217.89
$ws.Range("M5").Value = 7408.15
$ws.Range("N5").Value = 2068.81
$ws.Range("O5").Value = 47116.11
$ws.Range("Q5").Value = 1217.9

# --- Row 6: corrected data values ---
$ws.Range("A6").Value = "Dwell time (%)"
$ws.Range("B6").Value = 1.11
$ws.Range("C6").Value = 0.32
$ws.Range("D6").Value = 16.36
$ws.Range("E6").Value = 4.04
$ws.Range("F6").Value = 3.18
$ws.Range("G6").Value = 34.58
$ws.Range("H6").Value = 0.65
$ws.Range("I6").Value = 0.81
$ws.Range("J6").Value = 0.16
$ws.Range("K6").Value = 0.49
$ws.Range("L6").Value = 0.56
$ws.Range("M6").Value = 3.43
$ws.Range("N6").Value = 0.96
$ws.Range("O6").Value = 21.84
$ws.Range("Q6").Value = 0.56

# --- Row 7: corrected data values ---
$ws.Range("A7").Value = "Fixation duration (ms)"
$ws.Range("B7").Value = 340.8
$ws.Range("C7").Value = 350.29
$ws.Range("D7").Value = 289.36
$ws.Range("E7").Value = 322.53
$ws.Range("F7").Value = 285.7
$ws.Range("G7").Value = 252.75
$ws.Range("H7").Value = 280.31
$ws.Range("I7").Value = 350.34
$ws.Range("J7").Value = 350.36
$ws.Range("K7").Value = 355.83
$ws.Range("L7").Value = 304.47
$ws.Range("M7").Value = 352.77
$ws.Range("N7").Value = 413.76
$ws.Range("O7").Value = 395.93
$ws.Range("Q7").Value = 243.58

# --- Row 8: corrected data values ---
$ws.Range("A8").Value = "First fixation duration (ms)"
$ws.Range("B8").Value = 183.56
$ws.Range("C8").Value = 150.04
$ws.Range("D8").Value = 350.34
$ws.Range("E8").Value = 350.15
$ws.Range("F8").Value = 800.8
$ws.Range("G8").Value = 150.18
$ws.Range("H8").Value = 300.32
$ws.Range("I8").Value = 483.81
$ws.Range("J8").Value = 350.36
$ws.Range("K8").Value = 350.15
$ws.Range("L8").Value = 317
$ws.Range("M8").Value = 100.12
$ws.Range("N8").Value = 1201.31
$ws.Range("O8").Value = 550.56
$ws.Range("Q8").Value = 150.1
